$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 45454990
$ws.Cells.Item(33, 9).Value = 258.14285
$ws.Cells.Item(33, 11).Value = 258.14285
$ws.Cells.Item(33, 13).Value = -29.14285000000001

$ws.Cells.Item(53, 8).Value = 3188.9412
$ws.Cells.Item(53, 10).Value = 144.83333
$ws.Cells.Item(53, 12).Value = 144.83333
$ws.Cells.Item(53, 14).Value = -1418.83333

$ws.Cells.Item(106, 8).Value = 18443.059
$ws.Cells.Item(106, 10).Value = 1500
$ws.Cells.Item(106, 12).Value = 1500
$ws.Cells.Item(106, 14).Value = -2762

$ws.Cells.Item(111, 8).Value = 2005
$ws.Cells.Item(111, 9).Value = 1012.5
$ws.Cells.Item(111, 10).Value = 2997.5
$ws.Cells.Item(111, 11).Value = 3037.5
$ws.Cells.Item(111, 12).Value = 8992.5
$ws.Cells.Item(111, 13).Value = 29.5
$ws.Cells.Item(111, 14).Value = -15126.5

$ws.Cells.Item(116, 8).Value = 20000
$ws.Cells.Item(116, 9).Value = 20000
$ws.Cells.Item(116, 11).Value = 20000
$ws.Cells.Item(116, 13).Value = -16558

$ws.Cells.Item(132, 8).Value = 2030.7407
$ws.Cells.Item(132, 9).Value = 1921.28
$ws.Cells.Item(132, 11).Value = 5763.84
$ws.Cells.Item(132, 13).Value = -3233.84

$ws.Cells.Item(133, 8).Value = 83622.5
$ws.Cells.Item(133, 10).Value = 83622.5
$ws.Cells.Item(133, 12).Value = 83622.5
$ws.Cells.Item(133, 14).Value = -93742.5

$ws.Cells.Item(137, 8).Value = 5699
$ws.Cells.Item(137, 9).Value = 5874.375
$ws.Cells.Item(137, 10).Value = 4997.5
$ws.Cells.Item(137, 11).Value = 17623.125
$ws.Cells.Item(137, 12).Value = 14992.5
$ws.Cells.Item(137, 13).Value = -15073.125
$ws.Cells.Item(137, 14).Value = -20092.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 3271.0715
$ws.Cells.Item(2, 10).Value = 4042.1428
$ws.Cells.Item(2, 12).Value = 4042.1428
$ws.Cells.Item(2, 14).Value = -4268.1428

$ws.Cells.Item(32, 8).Value = 5383.061
$ws.Cells.Item(32, 9).Value = 1966.3243
$ws.Cells.Item(32, 10).Value = 15918
$ws.Cells.Item(32, 11).Value = 1966.3243
$ws.Cells.Item(32, 12).Value = 15918
$ws.Cells.Item(32, 13).Value = -1679.3243
$ws.Cells.Item(32, 14).Value = -16492

$ws.Cells.Item(45, 8).Value = 1882.7693
$ws.Cells.Item(45, 9).Value = 1882.7693
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 11).Value = 1882.7693
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 13).Value = -1505.7693
$ws.Cells.Item(45, 14).ClearContents()

$ws.Cells.Item(59, 8).Value = 0
$ws.Cells.Item(59, 10).Value = 0
$ws.Cells.Item(59, 12).Value = 0
$ws.Cells.Item(59, 14).ClearContents()

$ws.Cells.Item(74, 8).Value = 2754.946
$ws.Cells.Item(74, 9).Value = 2139.6875
$ws.Cells.Item(74, 11).Value = 2139.6875
$ws.Cells.Item(74, 13).Value = -1265.6875

$ws.Cells.Item(76, 8).Value = 21116.5
$ws.Cells.Item(76, 10).Value = 21116.5
$ws.Cells.Item(76, 12).Value = 21116.5
$ws.Cells.Item(76, 14).Value = -21792.5

$ws.Cells.Item(77, 8).Value = 2754.946
$ws.Cells.Item(77, 9).Value = 2139.6875
$ws.Cells.Item(77, 11).Value = 10698.4375
$ws.Cells.Item(77, 13).Value = -6330.4375

$ws.Cells.Item(79, 8).Value = 21116.5
$ws.Cells.Item(79, 10).Value = 21116.5
$ws.Cells.Item(79, 12).Value = 21116.5
$ws.Cells.Item(79, 14).Value = -23456.5

$ws.Cells.Item(109, 8).Value = 21377
$ws.Cells.Item(109, 10).Value = 21377
$ws.Cells.Item(109, 12).Value = 21377
$ws.Cells.Item(109, 14).Value = -24151

$ws.Cells.Item(110, 8).Value = 997.1667
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 14).ClearContents()

$ws.Cells.Item(116, 8).Value = 3271.0715
$ws.Cells.Item(116, 10).Value = 4042.1428
$ws.Cells.Item(116, 12).Value = 4042.1428
$ws.Cells.Item(116, 14).Value = -8630.1428

$ws.Cells.Item(123, 8).Value = 65028.5
$ws.Cells.Item(123, 10).Value = 65028.5
$ws.Cells.Item(123, 12).Value = 65028.5
$ws.Cells.Item(123, 14).Value = -74828.5

$ws.Cells.Item(132, 8).Value = 3801.3157
$ws.Cells.Item(132, 9).Value = 3904.3142
$ws.Cells.Item(132, 11).Value = 11712.9426
$ws.Cells.Item(132, 13).Value = -9182.942599999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 3271.0715
$ws.Cells.Item(3, 10).Value = 4042.1428
$ws.Cells.Item(3, 12).Value = 4042.1428
$ws.Cells.Item(3, 14).Value = -4270.1428

$ws.Cells.Item(61, 8).Value = 0
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 14).ClearContents()

$ws.Cells.Item(99, 8).Value = 55556930
$ws.Cells.Item(99, 9).Value = 62501332
$ws.Cells.Item(99, 11).Value = 62501332
$ws.Cells.Item(99, 13).Value = -62499834

$ws.Cells.Item(105, 8).Value = 1048.0435
$ws.Cells.Item(105, 9).Value = 1054.2727
$ws.Cells.Item(105, 11).Value = 1054.2727
$ws.Cells.Item(105, 13).Value = 692.7273

$ws.Cells.Item(132, 8).Value = 99995
$ws.Cells.Item(132, 10).Value = 99995
$ws.Cells.Item(132, 12).Value = 99995
$ws.Cells.Item(132, 14).Value = -110115

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2941.4167
$ws.Cells.Item(16, 9).Value = 2629.7
$ws.Cells.Item(16, 10).Value = 4500
$ws.Cells.Item(16, 11).Value = 2629.7
$ws.Cells.Item(16, 12).Value = 4500
$ws.Cells.Item(16, 13).Value = -2342.7
$ws.Cells.Item(16, 14).Value = -5074

$ws.Cells.Item(31, 8).Value = 2615.862
$ws.Cells.Item(31, 9).Value = 2384.1
$ws.Cells.Item(31, 11).Value = 2384.1
$ws.Cells.Item(31, 13).Value = -2089.1

$ws.Cells.Item(34, 8).Value = 2615.862
$ws.Cells.Item(34, 9).Value = 2384.1
$ws.Cells.Item(34, 11).Value = 2384.1
$ws.Cells.Item(34, 13).Value = -2182.1

$ws.Cells.Item(58, 8).Value = 6793
$ws.Cells.Item(58, 9).Value = 3451.25
$ws.Cells.Item(58, 11).Value = 3451.25
$ws.Cells.Item(58, 13).Value = -3248.25

$ws.Cells.Item(62, 8).Value = 4961.9
$ws.Cells.Item(62, 9).Value = 4402.222
$ws.Cells.Item(62, 11).Value = 4402.222
$ws.Cells.Item(62, 13).Value = -3778.222

$ws.Cells.Item(65, 8).Value = 4961.9
$ws.Cells.Item(65, 9).Value = 4402.222
$ws.Cells.Item(65, 11).Value = 22011.11
$ws.Cells.Item(65, 13).Value = -18891.11

$ws.Cells.Item(113, 8).Value = 2941.4167
$ws.Cells.Item(113, 9).Value = 2629.7
$ws.Cells.Item(113, 10).Value = 4500
$ws.Cells.Item(113, 11).Value = 2629.7
$ws.Cells.Item(113, 12).Value = 4500
$ws.Cells.Item(113, 13).Value = -459.6999999999998
$ws.Cells.Item(113, 14).Value = -8840

$ws.Cells.Item(134, 8).Value = 4990.6
$ws.Cells.Item(134, 10).Value = 4987
$ws.Cells.Item(134, 12).Value = 14961
$ws.Cells.Item(134, 14).Value = -20031

$ws.Cells.Item(136, 8).Value = 6793
$ws.Cells.Item(136, 9).Value = 3451.25
$ws.Cells.Item(136, 11).Value = 10353.75
$ws.Cells.Item(136, 13).Value = -7803.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 41709044
$ws.Cells.Item(4, 9).Value = 62563180
$ws.Cells.Item(4, 10).Value = 772.875
$ws.Cells.Item(4, 11).Value = 187689540
$ws.Cells.Item(4, 12).Value = 2318.625
$ws.Cells.Item(4, 13).Value = -187689428
$ws.Cells.Item(4, 14).Value = -2542.625

$ws.Cells.Item(33, 8).Value = 448.77777
$ws.Cells.Item(33, 9).Value = 108
$ws.Cells.Item(33, 10).Value = 619.1667
$ws.Cells.Item(33, 11).Value = 648
$ws.Cells.Item(33, 12).Value = 3715.0002
$ws.Cells.Item(33, 13).Value = -365
$ws.Cells.Item(33, 14).Value = -4281.0002

$ws.Cells.Item(46, 8).Value = 132939680
$ws.Cells.Item(46, 10).Value = 158734320
$ws.Cells.Item(46, 12).Value = 476202960
$ws.Cells.Item(46, 14).Value = -476203142

$ws.Cells.Item(80, 8).Value = 1975.3334
$ws.Cells.Item(80, 9).ClearContents()
$ws.Cells.Item(80, 10).Value = 1975.3334
$ws.Cells.Item(80, 11).ClearContents()
$ws.Cells.Item(80, 12).Value = 5926.0002
$ws.Cells.Item(80, 13).ClearContents()
$ws.Cells.Item(80, 14).Value = -7798.0002

$ws.Cells.Item(83, 8).Value = 1975.3334
$ws.Cells.Item(83, 9).ClearContents()
$ws.Cells.Item(83, 10).Value = 1975.3334
$ws.Cells.Item(83, 11).ClearContents()
$ws.Cells.Item(83, 12).Value = 17778.0006
$ws.Cells.Item(83, 13).ClearContents()
$ws.Cells.Item(83, 14).Value = -27138.0006

$ws.Cells.Item(92, 8).Value = 313.13635
$ws.Cells.Item(92, 9).Value = 267
$ws.Cells.Item(92, 10).Value = 323.3889
$ws.Cells.Item(92, 11).Value = 801
$ws.Cells.Item(92, 12).Value = 970.1667
$ws.Cells.Item(92, 13).Value = 447
$ws.Cells.Item(92, 14).Value = -3466.1667

$ws.Cells.Item(114, 8).Value = 14678.0625
$ws.Cells.Item(114, 9).Value = 209.375
$ws.Cells.Item(114, 11).Value = 628.125
$ws.Cells.Item(114, 13).Value = 2625.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 3299.2
$ws.Cells.Item(122, 9).Value = 3299.2
$ws.Cells.Item(122, 11).Value = 9897.599999999999
$ws.Cells.Item(122, 13).Value = -7447.599999999999

$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(123, 14).ClearContents()

$ws.Cells.Item(126, 8).Value = 5890.3125
$ws.Cells.Item(126, 9).Value = 4983.3335
$ws.Cells.Item(126, 11).Value = 14950.0005
$ws.Cells.Item(126, 13).Value = -12480.0005

$ws.Cells.Item(132, 8).Value = 4632.613
$ws.Cells.Item(132, 9).Value = 4411.6294
$ws.Cells.Item(132, 11).Value = 13234.8882
$ws.Cells.Item(132, 13).Value = -10704.8882

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1586.25
$ws.Cells.Item(46, 9).Value = 1088
$ws.Cells.Item(46, 10).Value = 2416.6667
$ws.Cells.Item(46, 11).Value = 1088
$ws.Cells.Item(46, 12).Value = 2416.6667
$ws.Cells.Item(46, 13).Value = -900
$ws.Cells.Item(46, 14).Value = -2792.6667

$ws.Cells.Item(55, 8).Value = 957.7692
$ws.Cells.Item(55, 10).Value = 1848.6
$ws.Cells.Item(55, 12).Value = 1848.6
$ws.Cells.Item(55, 14).Value = -2194.6

$ws.Cells.Item(136, 8).Value = 6332.5
$ws.Cells.Item(136, 9).Value = 5598
$ws.Cells.Item(136, 11).Value = 16794
$ws.Cells.Item(136, 13).Value = -14244

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1830.8125
$ws.Cells.Item(132, 9).Value = 1752.9231
$ws.Cells.Item(132, 11).Value = 5258.7693
$ws.Cells.Item(132, 13).Value = -2728.7693

$ws.Cells.Item(136, 8).Value = 4027.0435
$ws.Cells.Item(136, 9).Value = 4326.15
$ws.Cells.Item(136, 11).Value = 12978.45
$ws.Cells.Item(136, 13).Value = -10428.45
